$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 currently holds "Postsecondary" / its definition.
# Replace it with the first of the five new, more granular definitions,
# then insert four additional rows below it for the remaining four.

$ws.Range("A24").Value = "Apprenticeship or trades certificate or diploma"
$ws.Range("B24").Value = "Includes people who have obtained a certificate, diploma or equivalent in the skilled trades or in entry-level vocations, career and technical occupations. People who also have a postsecondary certificate, diploma or degree other than an apprenticeship or trades certificate or diploma are excluded from this category."

# Insert four new blank rows right after row 24 (rows 25-28), pushing the
# remainder of the table (Total household income, etc.) down to rows 29-39.
$ws.Range("25:28").Insert()

$ws.Range("A25").Value = "College, CEGEP or other non-university certificate or diploma"
$ws.Range("B25").Value = "Includes people who obtained a (non-trades) postsecondary certificate or diploma of 3 months or more from a public or private college, CEGEP, polytechnic, institute of technology, school of nursing, business school or vocational school. People who also have university certificates, diplomas or degrees are also excluded from this category."

$ws.Range("A26").Value = "University certificate or diploma below bachelor level"
$ws.Range("B26").Value = "Includes people who have obtained a university certificate or diploma below the bachelor level and who have not obtained any higher degrees, certificates or diplomas. University certificates or diplomas are commonly connected with professional associations in fields such as accounting, banking, insurance or public administration. The certificates and diplomas in this category do not require a bachelor's degree as a prerequisite."

$ws.Range("A27").Value = "Bachelor's degree or certificate/diploma above bachelor level"
$ws.Range("B27").Value = "Includes people who have obtained a bachelor's degree awarded by a degree-granting institution and people who have obtained a certificate or diploma that usually requires a bachelor's degree as a prerequisite. It excludes people who have obtained any higher degrees."

$ws.Range("A28").Value = "Graduate or professional degree"
$ws.Range("B28").Value = "Includes people who have obtained degrees in medicine, dentistry, veterinary medicine or optometry; people who have obtained a master's degree; and people who have obtained a doctorate degree awarded by a university."

# Match the final selection/active cell seen in the edited workbook.
$ws.Range("B23").Select()
